$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the NRIC value (B10) ---
$ws.Range("B10").Value = [char]0x202A + "S8927668C" + [char]0x202C
$ws.Range("B10").HorizontalAlignment = -4131

# --- Update email value in B2 and its hyperlink display text ---
$ws.Range("B2").Value = "sumit@circles.asia"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:sumit.testmail02@gmail.com", "", "", "sumit@circles.asia")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:sumit@circles.asia", "", "", "sumit@circles.asia")

$ws.Range("B2").Font.Color = 16711680
$ws.Range("B2").Font.Underline = $false
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10

$ws.Range("D2").Font.Color = 16711680
$ws.Range("D2").Font.Underline = $false
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Size = 10

# --- Update selection on sheet1 ---
$ws.Range("C17").Select()
